$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the two new sheets at the end: CodeSet, ActivityCode
# ---------------------------------------------------------------------------
$count = $wb.Worksheets.Count
$last = $wb.Worksheets.Item($count)

$codeSet = $wb.Worksheets.Add($null, $last)
$codeSet.Name = "CodeSet"

$activityCode = $wb.Worksheets.Add($null, $codeSet)
$activityCode.Name = "ActivityCode"

# ---------------------------------------------------------------------------
# 2. Populate CodeSet sheet
# ---------------------------------------------------------------------------
$codeSet.Range("A1").Value = "name"
$codeSet.Range("B1").Value = "slug"
$codeSet.Range("C1").Value = "description"
$codeSet.Range("A1:C1").Font.Bold = $true

$codeSet.Range("A2").Value = "Code set 1"
$codeSet.Range("B2").Value = "code_set1"

$codeSet.Range("A3").Value = "Code set 2"
$codeSet.Range("B3").Value = "code_set2"

$codeSet.Columns.Item(1).ColumnWidth = 35.29
$codeSet.Columns.Item(2).ColumnWidth = 24.17

$codeSet.Range("B2").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Populate ActivityCode sheet
# ---------------------------------------------------------------------------
$activityCode.Range("A1").Value = "codeset_slug"
$activityCode.Range("B1").Value = "activitycode"
$activityCode.Range("C1").Value = "label"
$activityCode.Range("D1").Value = "vertical_distribution_slug"
$activityCode.Range("A1:D1").Font.Bold = $true

$activityCode.Range("A2").Value = "code_set1"
$activityCode.Range("B2").Value = "1.3"
$activityCode.Range("C2").Value = "PublicPower"

$activityCode.Range("A3").Value = "code_set1"
$activityCode.Range("B3").Value = "1.A.4.b.i"
$activityCode.Range("C3").Value = "Industry"

$activityCode.Range("A4").Value = "code_set2"
$activityCode.Range("B4").Value = "A"
$activityCode.Range("C4").Value = "Combustion in the production and transformation of energy"

$activityCode.Range("B2:B5").NumberFormat = "@"

$activityCode.Columns.Item(1).ColumnWidth = 15.56
$activityCode.Columns.Item(2).ColumnWidth = 18.06
$activityCode.Columns.Item(3).ColumnWidth = 48.62

$activityCode.Range("C20").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Misc view tweaks on existing sheets
# ---------------------------------------------------------------------------
$pointSource = $wb.Worksheets.Item("PointSource")
$pointSource.Activate()
$pointSource.Range("G4").Select() | Out-Null

$activityCode.Activate()
